$d = $word.ActiveDocument

function Insert-XmlAtPoint($doc, $pos, $innerBodyXml) {
    $insPoint = $doc.Range($pos, $pos)
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$insPoint.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# Change 1: shorten the "These objectives will reveal..." sentence, dropping
# the second sentence about statistical procedures.
# ---------------------------------------------------------------------------
$oldPara1 = "These objectives will reveal important aspects of job satisfaction and help employers better understand what contributes to employee loyalty and satisfaction at work. The analysis process can utilize various statistical procedures, including Pearson's correlation coefficient, Spearman's rank correlation, Chi-square tests, etc., depending on the type and distribution of the data."
$newPara1 = "These objectives will reveal important aspects of job satisfaction and help employers better understand what contributes to employee loyalty and satisfaction at work. "
[void]$d.Content.Find.Execute($oldPara1, $true, $false, $false, $false, $false, $true, 1, $false, $newPara1, 2)

# ---------------------------------------------------------------------------
# Change 2: move <w:lastRenderedPageBreak/> from the "All roles that begin
# with Software..." run to the start of the following "Remaining entries
# classified..." run.
# ---------------------------------------------------------------------------

# 2a. Remove the page break from the "All roles that begin with..." run by
#     rebuilding that run's text without the break element.
$softwareText = 'All roles that begin with "Software..." were consolidated under the title "Software Engineer."'
$rng = $d.Content
[void]$rng.Find.Execute($softwareText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $rng.Start
$len = $softwareText.Length

$innerXml = '<w:body><w:p><w:r><w:t>' + $softwareText + '</w:t></w:r></w:p></w:body>'
Insert-XmlAtPoint $d $startPos $innerXml

$oldRng = $d.Range($startPos + $len, $startPos + $len + $len)
$oldRng.Text = ""

# 2b. Add the page break to the start of the "Remaining entries classified..."
#     run.
$remainingText = 'Remaining entries classified as "other" were uniformly relabeled as "Other."'
$rng2 = $d.Content
[void]$rng2.Find.Execute($remainingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos2 = $rng2.Start
$len2 = $remainingText.Length

$innerXml2 = '<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>' + $remainingText + '</w:t></w:r></w:p></w:body>'
Insert-XmlAtPoint $d $startPos2 $innerXml2

$oldRng2 = $d.Range($startPos2 + $len2, $startPos2 + $len2 + $len2)
$oldRng2.Text = ""

# ---------------------------------------------------------------------------
# Change 3: split the "Coefficient of Variation" run into two runs, moving
# <w:lastRenderedPageBreak/> to the start of the second run ("satisfaction
# with...").
# ---------------------------------------------------------------------------
$part1 = 'Coefficient of Variation (CV): The CV measures the ratio of the standard deviation to the mean and indicates the degree of data dispersion. A lower CV indicates less relative dispersion. In this case, '
$part2 = 'satisfaction with "Colleagues" has the lowest CV, suggesting that ratings in this area are more concentrated around the mean.'
$cvText = $part1 + $part2

$rng3 = $d.Content
[void]$rng3.Find.Execute($cvText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos3 = $rng3.Start
$len3 = $cvText.Length

$innerXml3 = '<w:body><w:p><w:r><w:t xml:space="preserve">' + $part1 + '</w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>' + $part2 + '</w:t></w:r></w:p></w:body>'
Insert-XmlAtPoint $d $startPos3 $innerXml3

$oldRng3 = $d.Range($startPos3 + $len3, $startPos3 + $len3 + $len3)
$oldRng3.Text = ""
